$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$apos = "'"  # used to force text-entry for numeric-looking Price values (mirrors typing an apostrophe in Excel)

$ws.Range("D2").Value = '71.399.53'
$ws.Range("E2").Value = '  -1.65%  '

$ws.Range("D3").Value = '3.963.39'
$ws.Range("E3").Value = '  -2.24%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = $apos + '543.83'
$ws.Range("E5").Value = '  +4.71%  '

$ws.Range("D6").Value = $apos + '149.17'
$ws.Range("E6").Value = '  +1.34%  '

$ws.Range("D7").Value = '3.954.06'
$ws.Range("E7").Value = '  -2.22%  '

$ws.Range("E8").Value = '  -6.78%  '

$ws.Range("E9").Value = '  +0.06%  '

$ws.Range("E10").Value = '  -4.41%  '

$ws.Range("D11").Value = $apos + '0.167'
$ws.Range("E11").Value = '  -5.88%  '

$ws.Range("D12").Value = $apos + '56.72'
$ws.Range("E12").Value = '  +18.48%  '

$ws.Range("D13").Value = $apos + '0.0000318'
$ws.Range("E13").Value = '  -2.96%  '

$ws.Range("D14").Value = $apos + '10.68'
$ws.Range("E14").Value = '  -4.69%  '

$ws.Range("D15").Value = '4.603.57'
$ws.Range("E15").Value = '  -2.23%  '

$ws.Range("D16").Value = '3.966.44'
$ws.Range("E16").Value = '  -2.46%  '

$ws.Range("D17").Value = $apos + '13.85'
$ws.Range("E17").Value = '  -2.27%  '

$ws.Range("D18").Value = $apos + '20.47'
$ws.Range("E18").Value = '  -3.99%  '

$ws.Range("D19").Value = $apos + '0.132'
$ws.Range("E19").Value = '  -1.13%  '

$ws.Range("E20").Value = '  -3.24%  '

$ws.Range("D21").Value = '71.241.18'
$ws.Range("E21").Value = '  -1.69%  '

$ws.Range("D22").Value = $apos + '425.67'
$ws.Range("E22").Value = '  -3.87%  '

$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = $apos + '97.29'
$ws.Range("E23").Value = '  -6.93%  '

$ws.Range("B24").Value = 'ImmutableX'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D24").Value = $apos + '3.58'
$ws.Range("E24").Value = '  -0.34%  '

$ws.Range("D25").Value = $apos + '4.24'
$ws.Range("E25").Value = '  +5.58%  '

$ws.Range("D26").Value = $apos + '14.34'
$ws.Range("E26").Value = '  -3.54%  '

$ws.Range("D27").Value = $apos + '11.56'
$ws.Range("E27").Value = '  +0.61%  '

$ws.Range("D28").Value = $apos + '10.73'
$ws.Range("E28").Value = '  -2.82%  '

$ws.Range("D29").Value = $apos + '3.77'
$ws.Range("E29").Value = '  +15.27%  '

$ws.Range("D30").Value = $apos + '5.89'
$ws.Range("E30").Value = '  +1.56%  '

$ws.Range("D31").Value = $apos + '36.59'
$ws.Range("E31").Value = '  -3.10%  '

$ws.Range("D32").Value = $apos + '7.81'
$ws.Range("E32").Value = '  +14.57%  '

$ws.Range("D33").Value = $apos + '51.52'
$ws.Range("E33").Value = '  +19.60%  '

$ws.Range("D34").Value = $apos + '693.85'
$ws.Range("E34").Value = '  +2.06%  '

$ws.Range("D35").Value = $apos + '13.40'
$ws.Range("E35").Value = '  -2.32%  '

$ws.Range("D36").Value = $apos + '0.131'
$ws.Range("E36").Value = '  -0.19%  '

$ws.Range("D37").Value = $apos + '64.96'
$ws.Range("E37").Value = '  -3.31%  '

$ws.Range("D38").Value = $apos + '0.436'
$ws.Range("E38").Value = '  +1.40%  '

$ws.Range("E39").Value = '  -0.56%  '

$ws.Range("D40").Value = '0.0₃0824'
$ws.Range("E40").Value = '  -4.79%  '

$ws.Range("D41").Value = $apos + '3.43'
$ws.Range("E41").Value = '  -3.03%  '

$ws.Range("E42").Value = '  +0.15%  '

$ws.Range("E43").Value = '  +0.18%  '

$ws.Range("D44").Value = $apos + '3.26'
$ws.Range("E44").Value = '  -0.10%  '

$ws.Range("D45").Value = $apos + '0.0483'
$ws.Range("E45").Value = '  -3.20%  '

$ws.Range("E46").Value = '  -5.92%  '

$ws.Range("D47").Value = $apos + '2.72'
$ws.Range("E47").Value = '  +0.32%  '

$ws.Range("D48").Value = $apos + '9.84'
$ws.Range("E48").Value = '  +7.95%  '

$ws.Range("D49").Value = $apos + '3.36'
$ws.Range("E49").Value = '  -4.12%  '

$ws.Range("D50").Value = $apos + '3.00'
$ws.Range("E50").Value = '  -2.07%  '

$ws.Range("E51").Value = '  +3.44%  '
